# stats.xlsx: correct a couple of figures in the "rainmaker" row on
# Sheet1 and leave the workbook focused on Sheet1 (instead of Sheet2) so
# both tables are ready to review/print, each with its last-used cell
# selection restored.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Corrected figures for row 4 ("rainmaker"): wins total and win pct.
$ws1.Range("B4").Value = 711
$ws1.Range("G4").Value = 0.754

# Preserve Sheet2's last selection before leaving it ...
$ws2.Activate()
[void]$ws2.Range("K6").Select()

# ... then make Sheet1 the active/selected tab with its own last
# selection restored.
$ws1.Activate()
[void]$ws1.Range("I11").Select()
